# Update the "取得日時" (acquired datetime) timestamps in the "ランサーズ" sheet
# from 2025-12-09 18:23:53 to 2025-12-09 18:33:17 for all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-09 18:23:53"
$newValue = "2025-12-09 18:33:17"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
